$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'287.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.95%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'29.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'2.28%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.095"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.08%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.06688"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'3.16%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'7.344"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.70%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.408"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'1.14%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.372"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'3.21%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9194"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.72%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1592"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'3.10%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.06824"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'5.07%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07607"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-0.12%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.02932"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-1.77%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.08978"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.25%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001573"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-1.77%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.04505"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.01%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.0006463"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-1.02%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.006306"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'3.08%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'3.450"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Value = "'0.3213"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'1.10%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'-2.38%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'4.069"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'2.39%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.1583"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'1.76%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.001190"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.65%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004110"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-4.87%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0001199"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'1.54%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.0001618"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'-1.16%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").Value = "'0.04266"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'2.88%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006742"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'0.61%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1237"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.002229"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'5.52%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'12.70%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005703"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'5.74%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'-3.59%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.01307"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-29.38%"
$ws.Range("E47").Style = "Normal"
